{"js": "// The self-assessment paragraph originally rambled into naming specific\n// teammates (\"Jacob\", \"Spencer\") and their contributions. The edit trims\n// that tangent and replaces it with a short transition sentence, while\n// keeping the surrounding sentences (the \"In order to work around this...\"\n// lead-in and the closing \"...I was successful in implementing my desired\n// features...\" sentence) intact.\n//\n// Rather than depending on exact pre-existing run boundaries (which are an\n// implementation detail of how the document happened to be saved), locate\n// the target text by searching body text and replace it in place. This\n// mirrors what a human editing in Word (or a Find & Replace macro) would\n// do, and is robust to minor run-splitting differences in the source file.\n\nconst body = context.document.body;\n\nconst oldText =\n  \" found myself doing most of the work\\u2026if I can call that an obstacle. \" +\n  \"If any team members deserve recognition, I would say Jacob did the most \" +\n  \"work besides me. Spencer did very little work on the programming and on \" +\n  \"the assignments overall; but, he did get the poster printed off which is \" +\n  \"one thing. I also will add that Spencer was very good at talking to \" +\n  \"people at the Expo. Anyways, to get back to the point of this paragraph\";\n\nconst newText = \"n addition to this\";\n\nconst searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target passage not found; document may already be edited.\");\n}\n\nsearchResults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The self-assessment paragraph originally rambled into naming specific\n# teammates (\"Jacob\", \"Spencer\") and their contributions. This edit trims\n# that tangent and replaces it with a short transition sentence, while\n# keeping the surrounding sentences (the \"In order to work around this...\"\n# lead-in and the closing \"...I was successful in implementing my desired\n# features...\" sentence) intact.\n\n$d = $word.ActiveDocument\n\n# U+2026 HORIZONTAL ELLIPSIS, built from its code point to avoid any source\n# encoding ambiguity.\n$ellipsis = [char]0x2026\n\n$findText = \" found myself doing most of the work\" + $ellipsis + `\n    \"if I can call that an obstacle. If any team members deserve recognition, \" + `\n    \"I would say Jacob did the most work besides me. Spencer did very little \" + `\n    \"work on the programming and on the assignments overall; but, he did get \" + `\n    \"the poster printed off which is one thing. I also will add that Spencer \" + `\n    \"was very good at talking to people at the Expo. Anyways, to get back to \" + `\n    \"the point of this paragraph\"\n\n$replaceText = \"n addition to this\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne)\n\nif (-not $found) {\n    throw \"Target passage not found; document may already be edited.\"\n}\n"}
